$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 41-42, shifting existing rows 41-100 down to 43-102
$ws.Rows("41:42").Insert()

# Fill new row 41
$ws.Range("A41").Value = 7
$ws.Range("B41").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C41").Value = "Ñuble"
$ws.Range("D41").Value = 44413
$ws.Range("E41").Value = 16
$ws.Range("F41").Value = 100112002
$ws.Range("G41").Value = "Pimiento"
$ws.Range("H41").Value = "Zafiro rojo"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 160
$ws.Range("K41").Value = 24000
$ws.Range("L41").Value = 25000
$ws.Range("M41").Value = 24500
$ws.Range("N41").Value = "`$/caja 15 kilos"
$ws.Range("O41").Value = "Región de Arica y Parinacota"
$ws.Range("P41").Value = 1633
$ws.Range("Q41").Value = 15
$ws.Range("R41").Value = "Hortaliza"

# Fill new row 42
$ws.Range("A42").Value = 7
$ws.Range("B42").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C42").Value = "Ñuble"
$ws.Range("D42").Value = 44413
$ws.Range("E42").Value = 16
$ws.Range("F42").Value = 100112002
$ws.Range("G42").Value = "Pimiento"
$ws.Range("H42").Value = "Zafiro verde"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 120
$ws.Range("K42").Value = 21000
$ws.Range("L42").Value = 22000
$ws.Range("M42").Value = 21500
$ws.Range("N42").Value = "`$/caja 15 kilos"
$ws.Range("O42").Value = "Región de Arica y Parinacota"
$ws.Range("P42").Value = 1433
$ws.Range("Q42").Value = 15
$ws.Range("R42").Value = "Hortaliza"
